# lightsheet-metadata.xlsx: move changes to v2; regen
#
# - resolution_z_value/resolution_z_unit columns become range_z_value/range_z_unit
# - new step_z_value, increment_z_value, increment_z_unit columns are inserted
#   right after them (shifting number_of_antibodies .. data_path three columns right)
# - "resolution_z_unit list" sheet is renamed "range_z_unit list"
# - a new "increment_z_unit list" sheet is added (same nm/um options)
# - schema version bumped from 1 to 2

$wb = $excel.ActiveWorkbook

function Set-CellComment($rng, [string]$text) {
    if ($rng.Comment -ne $null) {
        $rng.Comment.Delete()
    }
    $rng.AddComment($text) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Rename the resolution_z_unit list sheet, add the increment_z_unit list
# ---------------------------------------------------------------------------

$wsRangeUnit = $wb.Worksheets.Item("resolution_z_unit list")
$wsRangeUnit.Name = "range_z_unit list"

$wsIncUnit = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsRangeUnit)
$wsIncUnit.Name = "increment_z_unit list"
$wsIncUnit.Range("A1").Value = "nm"
$wsIncUnit.Range("A2").Value = "um"

# ---------------------------------------------------------------------------
# 2. Bump the schema version list value from 1 to 2
# ---------------------------------------------------------------------------

$wsVersion = $wb.Worksheets.Item("version list")
$wsVersion.Range("A1").Value = "'2"

# ---------------------------------------------------------------------------
# 3. Main "Export as TSV" sheet: shift columns W:AA right by three (to Z:AD),
#    rename U/V, and populate the three new columns W/X/Y.
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Export as TSV")

# Work right-to-left so we never overwrite a value before we've copied it.
$ws.Range("AD1").Value = "data_path"
$ws.Range("AC1").Value = "contributors_path"
$ws.Range("AB1").Value = "antibodies_path"
$ws.Range("AA1").Value = "number_of_channels"
$ws.Range("Z1").Value  = "number_of_antibodies"

$ws.Range("U1").Value = "range_z_value"
$ws.Range("V1").Value = "range_z_unit"
$ws.Range("W1").Value = "step_z_value"
$ws.Range("X1").Value = "increment_z_value"
$ws.Range("Y1").Value = "increment_z_unit"

# New header cells need the bold/centered/wrap-text header style; copy it
# from an existing header cell. (AB1:AD1 are brand-new cells too, even
# though they hold "old" header text, so they need the style re-applied.)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("W1:Y1").PasteSpecial(-4122) | Out-Null
$ws.Range("AB1:AD1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- comments -----------------------------------------------------------

Set-CellComment $ws.Range("U1")  "The total range of the z axis."
Set-CellComment $ws.Range("V1")  "The unit of range_z_value."
Set-CellComment $ws.Range("W1")  "The number of optical sections in z axis range."
Set-CellComment $ws.Range("X1")  "The distance between sequential optical sections."
Set-CellComment $ws.Range("Y1")  "The units of increment z value."
Set-CellComment $ws.Range("Z1")  "Number of antibodies"
Set-CellComment $ws.Range("AA1") "Number of fluorescent channels imaged during each cycle."
Set-CellComment $ws.Range("AB1") "Relative path to file with antibody information for this dataset."
Set-CellComment $ws.Range("AC1") "Relative path to file with ORCID IDs for contributors for this dataset."
Set-CellComment $ws.Range("AD1") "Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions."

# --- data validation ------------------------------------------------------

# version list error message now references "2"
$vA = $ws.Range("A2:A1048576").Validation
$vA.Modify(3, 1, 1, "'version list'!`$A`$1:`$A`$1")
$vA.ErrorTitle = "Value must come from list"
$vA.ErrorMessage = "Value must be one of: 2."

# range_z_unit list source renamed
$vV = $ws.Range("V2:V1048576").Validation
$vV.Modify(3, 1, 1, "'range_z_unit list'!`$A`$1:`$A`$2")
$vV.ErrorTitle = "Value must come from list"
$vV.ErrorMessage = "Value must be one of: nm / um."

# the old W validation (number_of_antibodies, whole) is repurposed in place
# for the new step_z_value column (decimal)
$vW = $ws.Range("W2:W1048576").Validation
$vW.Modify(2, 1, 1, "-1e+307", "1e+307")
$vW.ErrorTitle = "Not a number"
$vW.ErrorMessage = "The values in this column must be numbers."

# the old X validation (number_of_channels, whole) is repurposed in place
# for the new increment_z_value column (decimal)
$vX = $ws.Range("X2:X1048576").Validation
$vX.Modify(2, 1, 1, "-1e+307", "1e+307")
$vX.ErrorTitle = "Not a number"
$vX.ErrorMessage = "The values in this column must be numbers."

# new Y = increment_z_unit, list
$vY = $ws.Range("Y2:Y1048576").Validation
$vY.Add(3, 1, 1, "'increment_z_unit list'!`$A`$1:`$A`$2")
$vY.ErrorTitle = "Value must come from list"
$vY.ErrorMessage = "Value must be one of: nm / um."

# number_of_antibodies validation (whole) moves from old W to new Z
$vZ = $ws.Range("Z2:Z1048576").Validation
$vZ.Add(1, 1, 1, "-2147483647", "2147483647")
$vZ.ErrorTitle = "Not an integer"
$vZ.ErrorMessage = "The values in this column must be integers."

# number_of_channels validation (whole) moves from old X to new AA
$vAA = $ws.Range("AA2:AA1048576").Validation
$vAA.Add(1, 1, 1, "-2147483647", "2147483647")
$vAA.ErrorTitle = "Not an integer"
$vAA.ErrorMessage = "The values in this column must be integers."
